$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E74) listed the periods in descending
# order (2105 .. 1607). The sheet is updated so the periods are listed in
# ascending order (1607 .. 2105) instead - i.e. the 59 data rows are
# reversed. The "Valor Mora" (F) value of 48000 stays attached to period
# 2105 (all other periods carry 60000), so after the reversal that value
# moves from the first data row (16) to the last (74).

$year = 16
$month = 7
for ($row = 16; $row -le 74; $row++) {
    $period = "{0:D2}{1:D2}" -f $year, $month
    $ws.Range("E$row").Value = $period

    $month = $month + 1
    if ($month -gt 12) {
        $month = 1
        $year = $year + 1
    }
}

$ws.Range("F16").Value = 60000
$ws.Range("F74").Value = 48000
